# Applies the "Add files via upload" revision:
#  - GeneralVariables:  clear all example values in column B (rows 2-16)
#  - PipetteVariables:  clear all example values in column B (rows 2-8)
#  - PerPlateVariables: rename "Plate Example1"/"Plate Example2" headers to
#    "Plate 1"/"Plate 2", drop the "Plate Example3" column entirely and
#    clear the remaining example data
#  - refresh the selections left behind on each sheet

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("GeneralVariables")
$wsPipette = $wb.Worksheets.Item("PipetteVariables")
$wsPlate   = $wb.Worksheets.Item("PerPlateVariables")

# --- GeneralVariables: blank out the sample Value column -------------------
$wsGeneral.Range("B2:B16").ClearContents()

# --- PipetteVariables: blank out the sample Value column -------------------
$wsPipette.Range("B2:B8").ClearContents()

# --- PerPlateVariables: drop the 3rd plate column, rename the others, and --
# --- clear the leftover sample data ----------------------------------------
$tbl = $wsPlate.ListObjects.Item("Tabla2")

$wsPlate.Range("B2:D7").ClearContents()
$tbl.ListColumns.Item(4).Delete()
$wsPlate.Columns("D").Delete()
$wsPlate.Columns("C").Delete()

$wsPlate.Range("B1").Value = "Plate 1"
$wsPlate.Range("C1").Value = "Plate 2"

# --- restore the on-screen selections for each sheet ------------------------
[void]$wsPlate.Range("C2").Select()
[void]$wsPipette.Range("B2:B8").Select()
[void]$wsGeneral.Range("B2:B16").Select()
[void]$wsGeneral.Activate()
